# Apply the MIMAG template update:
#  - bump the template Version value on the isa_template sheet (1.0.0 -> 1.0.1)
#  - rename a few Assay table columns (header row) to the new naming scheme
#  - add one example data row underneath the Assay table headers

$wb = $excel.ActiveWorkbook

# ---- isa_template sheet: bump Version ----
$ws1 = $wb.Worksheets.Item("isa_template")
$ws1.Range("B4").Value = "1.0.1"

# ---- Assay sheet: rename headers ----
$ws2 = $wb.Worksheets.Item("Assay")
$ws2.Range("A1").Value  = "Input [Data]"
$ws2.Range("H1").Value  = "Parameter [sequence assembly algorithm version]"
$ws2.Range("I1").Value  = "Term Source REF (DPBO:0000060)"
$ws2.Range("J1").Value  = "Term Accession Number (DPBO:0000060)"
$ws2.Range("AI1").Value = "Output [Data]"

# ---- Assay sheet: add example data row to the table ----
$tbl = $ws2.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

$ws2.Range("B2").Value  = "Forest soil metagenome"
$ws2.Range("E2").Value  = "metaSPAdes"
$ws2.Range("H2").Value  = "3.11.0"
$ws2.Range("K2").Value  = "kmer set 21,33,55,77,99,121, default parameters otherwise"
$ws2.Range("N2").Value  = "metabat"
$ws2.Range("Q2").Value  = "homology search, kmer, coverage, codon usage, combination"
$ws2.Range("T2").Value  = "high"
$ws2.Range("W2").Value  = "Bacterial 16S RNA"
$ws2.Range("X2").Value  = "NCIT"
$ws2.Range("Y2").Value  = "https://bioregistry.io/NCIT:C105370"
$ws2.Range("Z2").Value  = "checkm"

# completeness/contamination scores are plain text percentages in the
# source data (not numeric Percentage-formatted cells), so force text
# interpretation before assigning - otherwise Excel auto-converts "90%"
# into the number 0.9 with a Percentage number format.
$ws2.Range("AC2").NumberFormat = "@"
$ws2.Range("AC2").Value = "90%"
$ws2.Range("AF2").NumberFormat = "@"
$ws2.Range("AF2").Value = "4%"
